# Update the "Datos actualizados" timestamp text in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 12:22"

# Row 9: Valencia/Valencia - update counts
$ws.Range("B9").Value = 4188
$ws.Range("C9").Value = 946
$ws.Range("D9").Value = 2881
$ws.Range("E9").Value = 361

# Row 12: La Rioja - update counts
$ws.Range("B12").Value = 3026
$ws.Range("C12").Value = 1172
$ws.Range("D12").Value = 1672
$ws.Range("E12").Value = 182

# Row 13: was Araba/Alava, now Alacant/Alicante - update name and counts
$ws.Range("A13").Value = "Alacant/Alicante"
$ws.Range("B13").Value = 2850
$ws.Range("C13").Value = 632
$ws.Range("D13").Value = 1908
$ws.Range("E13").Value = 310

# Row 14: was Alacant/Alicante, now Araba/Alava - update name and counts
$ws.Range("A14").Value = "Araba/Alava"
$ws.Range("B14").Value = 2806
$ws.Range("C14").Value = 4151
$ws.Range("D14").Value = 4663
$ws.Range("E14").Value = 229

# Row 25: was Granada, now Cantabria - update name and counts
$ws.Range("A25").Value = "Cantabria"
$ws.Range("B25").Value = 1619
$ws.Range("C25").Value = 214
$ws.Range("D25").Value = 1307
$ws.Range("E25").Value = 98

# Row 26: was Cantabria, now Granada - update name and counts
$ws.Range("A26").Value = "Granada"
$ws.Range("B26").Value = 1600
$ws.Range("C26").Value = 177
$ws.Range("D26").Value = 1285
$ws.Range("E26").Value = 138

# Row 37: was Aragon, now Castello/Castellon - update name and counts
$ws.Range("A37").Value = "Castello/Castellon"
$ws.Range("B37").Value = 921
$ws.Range("C37").Value = 194
$ws.Range("D37").Value = 631
$ws.Range("E37").Value = 96

# Row 38: was Cadiz, now Aragon - update name and counts
$ws.Range("A38").Value = "Aragon"
$ws.Range("B38").Value = 907
$ws.Range("C38").Value = 29
$ws.Range("D38").Value = 838
$ws.Range("E38").Value = 40

# Row 39: was Castello/Castellon, now Cadiz - update name and counts
$ws.Range("A39").Value = "Cadiz"
$ws.Range("B39").Value = 901
$ws.Range("C39").Value = 118
$ws.Range("D39").Value = 747
$ws.Range("E39").Value = 36
